$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text must be pre-formatted as Text so Excel
# stores the literal string (matching the source inlineStr cells) instead of
# silently converting "590.98" etc. into a number.
$textCells = @("D5","D6","D8","D10","D11","D12","D13","D14","D16","D19","D20","D21","D22","D23","D24","D26","D27","D29","D30","D31","D32","D33","D34","D37","D38","D39","D41","D42","D43","D44","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed coin data (price + 1h volume deltas, with a handful of
# row re-labels where the scrape re-ordered near-duplicate entries).
$ws.Range('D2').Value = '67.387.01'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '2.552.35'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '590.98'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = '174.52'
$ws.Range('E6').Value = '  +5.05%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.529'
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').Value = '2.550.09'
$ws.Range('E9').Value = '  -2.42%  '
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D11').Value = '0.162'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '0.353'
$ws.Range('E12').Value = '  -2.67%  '
$ws.Range('D13').Value = '5.17'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').Value = '3.012.18'
$ws.Range('E15').Value = '  -2.49%  '
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').Value = '67.210.29'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '2.558.18'
$ws.Range('E18').Value = '  -3.61%  '
$ws.Range('D19').Value = '8.13'
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').Value = '11.43'
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('D21').Value = '356.31'
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').Value = '4.70'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '70.18'
$ws.Range('E26').Value = '  +1.69%  '
$ws.Range('D27').Value = '9.89'
$ws.Range('E27').Value = '  -6.29%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.671.86'
$ws.Range('E28').Value = '  -3.01%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '0.0000100'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '538.58'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '8.30'
$ws.Range('E32').Value = '  +5.05%  '
$ws.Range('D33').Value = '1.35'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').Value = '158.01'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').Value = '18.83'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').Value = '0.358'
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('D42').Value = '1.82'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').Value = '2.56'
$ws.Range('E44').Value = '  +6.58%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '39.79'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').Value = '152.69'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').Value = '0.565'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0282'
$ws.Range('E49').Value = '  -6.57%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '3.73'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = '1.74'
$ws.Range('E51').Value = '  +2.09%  '

# Restore the default cell style on the text-forced cells (keeps formatting
# identical to the untouched cells; the stored value stays text because the
# content was already committed as a string above).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
